$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Y"
$ws.Range("H4").Value = "Y"
$ws.Range("H6").Value = "Y"
$ws.Range("H8").Value = "Y"
$ws.Range("H9").Value = "Y"
$ws.Range("H13").Value = "Y"
$ws.Range("H16").Value = "Y"
$ws.Range("H27").Value = "Y"

$ws.Range("H22").Select()
